$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B9: extend the sum formula with an extra term (5+1+2 -> 5+1+2+3)
$ws.Range("B9").Formula = "=5+1+2+3"

# B10: fill in a new time entry
$ws.Range("B10").Value = 6.25

# Move the active selection to B11 (matches post-edit cursor position)
$ws.Range("B11").Select()
